$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2564746666666666
$ws.Range("H2").Value = 0.7694239999999999
$ws.Range("I2").Value = 0.1818007399394835
$ws.Range("J2").Value = 0.1818007399394835
$ws.Range("M2").Value = 0.029424
$ws.Range("N2").Value = 0.08827199999999999
$ws.Range("O2").Value = 0.1473063425232919
$ws.Range("P2").Value = 0.1473063425232919
$ws.Range("Q2").Value = 0.007546510591999998
$ws.Range("R2").Value = 0.06791859532799999
$ws.Range("S2").Value = 0.02678040206851348
$ws.Range("T2").Value = 0.02678040206851348

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2564746666666666
$ws.Range("H3").Value = 0.7694239999999999
$ws.Range("I3").Value = 0.1818007399394835
$ws.Range("J3").Value = 0.1818007399394835
$ws.Range("O3").Value = 0.852693657476708
$ws.Range("P3").Value = 0.852693657476708
$ws.Range("Q3").Value = 0.04368353465066666
$ws.Range("R3").Value = 0.393151811856
$ws.Range("S3").Value = 0.15502033787097
$ws.Range("T3").Value = 0.1550203378709701

$ws.Range("I4").Value = 0.7694380609030022
$ws.Range("J4").Value = 0.7694380609030022
$ws.Range("M4").Value = 0.029424
$ws.Range("N4").Value = 0.08827199999999999
$ws.Range("O4").Value = 0.1473063425232919
$ws.Range("P4").Value = 0.1473063425232919
$ws.Range("Q4").Value = 0.03193921255999999
$ws.Range("R4").Value = 0.2874529130399999
$ws.Range("S4").Value = 0.1133431065498352
$ws.Range("T4").Value = 0.1133431065498352

$ws.Range("I5").Value = 0.7694380609030022
$ws.Range("J5").Value = 0.7694380609030022
$ws.Range("O5").Value = 0.852693657476708
$ws.Range("P5").Value = 0.852693657476708
$ws.Range("S5").Value = 0.6560949543531669
$ws.Range("T5").Value = 0.6560949543531669

$ws.Range("G6").Value = 0.06878966666666667
$ws.Range("I6").Value = 0.0487611991575143
$ws.Range("J6").Value = 0.0487611991575143
$ws.Range("M6").Value = 0.029424
$ws.Range("N6").Value = 0.08827199999999999
$ws.Range("O6").Value = 0.1473063425232919
$ws.Range("P6").Value = 0.1473063425232919
$ws.Range("Q6").Value = 0.002024067151999999
$ws.Range("R6").Value = 0.018216604368
$ws.Range("S6").Value = 0.007182833904943256
$ws.Range("T6").Value = 0.007182833904943256

$ws.Range("G7").Value = 0.06878966666666667
$ws.Range("I7").Value = 0.0487611991575143
$ws.Range("J7").Value = 0.0487611991575143
$ws.Range("O7").Value = 0.852693657476708
$ws.Range("P7").Value = 0.852693657476708
$ws.Range("Q7").Value = 0.01171646239566667
$ws.Range("S7").Value = 0.04157836525257104
$ws.Range("T7").Value = 0.04157836525257104
